# Update the cryptocurrency price/volume table with refreshed values.
# Column D (Price) and E (Volume(1h)) are stored as plain text in the
# workbook, so for any D-column value that looks like a pure number we
# force the cell to Text format before assigning, then restore the
# "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.957.39"
$ws.Range("E2").Value = "  -1.93%  "

$ws.Range("D3").Value = "2.449.85"
$ws.Range("E3").Value = "  -3.67%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.26%  "

$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("E11").Value = "  -4.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.07%  "

$ws.Range("D13").Value = "2.883.44"
$ws.Range("E13").Value = "  -3.59%  "

$ws.Range("D14").Value = "57.886.89"
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000133"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("D17").Value = "2.451.44"
$ws.Range("E17").Value = "  -3.46%  "

$ws.Range("E18").Value = "  -3.58%  "

$ws.Range("E19").Value = "  -2.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "313.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.47%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.404"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.74%  "

$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("D26").Value = "2.566.40"
$ws.Range("E26").Value = "  -2.97%  "

$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("E28").Value = "  -2.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "175.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.73%  "

$ws.Range("E30").Value = "  -3.20%  "

$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("E33").Value = "  -6.86%  "

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("E37").Value = "  -7.26%  "

$ws.Range("E38").Value = "  -4.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.21%  "

$ws.Range("E40").Value = "  +2.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.94%  "

$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.584"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.54%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "258.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.66%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.75%  "

$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0491"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.26%  "

$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.14%  "
